# Add daily power records to the comforter-cda sheet/table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the individual D/E/F formulas in rows 34-41 into shared
# formula groups (matches how Excel re-shares a formula range after an
# edit touches the whole block), and update the End Time for 9/22 (row 41).
$ws.Range("D34:D41").Formula = "=(C34-B34)* 1440"
$ws.Range("E34:E41").Formula = "=IF(C34>B34, (C34-B34)*1440, (B34-C34)*1440)"
$ws.Range("F34:F41").Formula = "=ABS((C34-B34)*1440)"
$ws.Range("C41").Value = 0.94027777777777777

# --- Grow the table by two rows (9/23 and 9/24 power records) so the
# table ref/autoFilter/dimension all expand to F43.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 42: 2018-09-23, power off at midnight, back on at 08:16:00.
$ws.Range("A42").Value = 43366
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 0.3444444444444445
$ws.Range("D42").Formula = "=(C42-B42)* 1440"
$ws.Range("E42").Formula = "=IF(C42>B42, (C42-B42)*1440, (B42-C42)*1440)"
$ws.Range("F42").Formula = "=ABS((C42-B42)*1440)"

# Row 43: 2018-09-24, no outage recorded yet (Start/End Time blank).
$ws.Range("A43").Value = 43367
$ws.Range("D43").Formula = "=(C43-B43)* 1440"
$ws.Range("E43").Formula = "=IF(C43>B43, (C43-B43)*1440, (B43-C43)*1440)"
$ws.Range("F43").Formula = "=ABS((C43-B43)*1440)"

# Leave the selection on the newly-added row, like the author did.
$ws.Range("B43").Select() | Out-Null
